{"js": "// Update the worksheet date and all two-digit-by-two-digit multiplication\n// answers to the newly generated problem set.\nconst replacements = [\n  [\"2024-07-31 Wednesday\", \"2024-08-01 Thursday\"],\n  [\"92\u00d787=8004\", \"97\u00d713=1261\"],\n  [\"16\u00d730=480\", \"70\u00d746=3220\"],\n  [\"15\u00d759=885\", \"29\u00d777=2233\"],\n  [\"38\u00d763=2394\", \"43\u00d766=2838\"],\n  [\"47\u00d729=1363\", \"12\u00d754=648\"],\n  [\"60\u00d758=3480\", \"18\u00d784=1512\"],\n  [\"96\u00d777=7392\", \"17\u00d776=1292\"],\n  [\"65\u00d774=4810\", \"58\u00d752=3016\"],\n  [\"79\u00d728=2212\", \"20\u00d742=840\"],\n  [\"39\u00d758=2262\", \"95\u00d795=9025\"],\n  [\"61\u00d738=2318\", \"54\u00d779=4266\"],\n  [\"32\u00d768=2176\", \"60\u00d737=2220\"],\n  [\"32\u00d723=736\", \"99\u00d712=1188\"],\n  [\"18\u00d779=1422\", \"56\u00d727=1512\"],\n  [\"13\u00d755=715\", \"55\u00d739=2145\"],\n  [\"24\u00d741=984\", \"81\u00d714=1134\"],\n  [\"86\u00d766=5676\", \"62\u00d799=6138\"],\n  [\"66\u00d778=5148\", \"53\u00d770=3710\"],\n  [\"65\u00d721=1365\", \"40\u00d787=3480\"],\n  [\"36\u00d769=2484\", \"39\u00d770=2730\"],\n  [\"60\u00d721=1260\", \"64\u00d793=5952\"],\n  [\"80\u00d740=3200\", \"15\u00d786=1290\"],\n  [\"15\u00d765=975\", \"99\u00d717=1683\"],\n  [\"17\u00d747=799\", \"94\u00d762=5828\"],\n  [\"27\u00d734=918\", \"64\u00d723=1472\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: true,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and all two-digit-by-two-digit multiplication\n# answers to the newly generated problem set.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2024-07-31 Wednesday\", \"2024-08-01 Thursday\"),\n  @(\"92\u00d787=8004\", \"97\u00d713=1261\"),\n  @(\"16\u00d730=480\", \"70\u00d746=3220\"),\n  @(\"15\u00d759=885\", \"29\u00d777=2233\"),\n  @(\"38\u00d763=2394\", \"43\u00d766=2838\"),\n  @(\"47\u00d729=1363\", \"12\u00d754=648\"),\n  @(\"60\u00d758=3480\", \"18\u00d784=1512\"),\n  @(\"96\u00d777=7392\", \"17\u00d776=1292\"),\n  @(\"65\u00d774=4810\", \"58\u00d752=3016\"),\n  @(\"79\u00d728=2212\", \"20\u00d742=840\"),\n  @(\"39\u00d758=2262\", \"95\u00d795=9025\"),\n  @(\"61\u00d738=2318\", \"54\u00d779=4266\"),\n  @(\"32\u00d768=2176\", \"60\u00d737=2220\"),\n  @(\"32\u00d723=736\", \"99\u00d712=1188\"),\n  @(\"18\u00d779=1422\", \"56\u00d727=1512\"),\n  @(\"13\u00d755=715\", \"55\u00d739=2145\"),\n  @(\"24\u00d741=984\", \"81\u00d714=1134\"),\n  @(\"86\u00d766=5676\", \"62\u00d799=6138\"),\n  @(\"66\u00d778=5148\", \"53\u00d770=3710\"),\n  @(\"65\u00d721=1365\", \"40\u00d787=3480\"),\n  @(\"36\u00d769=2484\", \"39\u00d770=2730\"),\n  @(\"60\u00d721=1260\", \"64\u00d793=5952\"),\n  @(\"80\u00d740=3200\", \"15\u00d786=1290\"),\n  @(\"15\u00d765=975\", \"99\u00d717=1683\"),\n  @(\"17\u00d747=799\", \"94\u00d762=5828\"),\n  @(\"27\u00d734=918\", \"64\u00d723=1472\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n  $find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
